$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selected/active cell on the sheet view
$ws.Range("W15").Select()

# New values for rows 4-10, columns L..X (the data that changed per the diff)
$data = @{
    4  = @{ L=1; M=1; N=3; O=4; P=7; Q=7; R=7; S=6; T=6; U=6; V=6; W=4; X=3 }
    5  = @{ M=2; N=4; O=4; P=8; Q=8; R=8; S=8; T=6; X=4 }
    6  = @{ L=2; M=2; N=4; O=4; P=8; Q=8; R=8; S=8; T=6; U=6; V=4; W=4; X=4 }
    7  = @{ L=2; M=2; N=4; O=4; Q=8; R=8; T=6; U=6; X=4 }
    8  = @{ L=2; M=2; N=4; O=5; P=8; Q=8; R=8; S=8; T=6; U=6; V=6; W=4; X=3 }
    9  = @{ N=5; O=6; P=9; Q=9; R=9; S=8; T=8; U=6; V=6; W=4; X=3 }
    10 = @{ M=3; O=6; P=9; Q=9; R=9; S=8; T=8; U=6; V=6; W=4; X=3 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $cellRef = "$col$row"
        $ws.Range($cellRef).Value = $cols[$col]
    }
}
